$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a single date value, repeated down every
# data row. The sheet was refreshed, bumping that date by one day
# (2023-09-09 -> 2023-09-10, serial 45178 -> 45179) for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 514 }

$ws.Range("C2:C$lastRow").Value = 45179
